# Applies the "demand" sheet update:
#  - Flip visibility of rows filtered by Scenario (National Trends rows get
#    hidden, Distributed Energy rows get shown) to match the new AutoFilter
#    criteria.
#  - Append a new H2 demand data row for node NOS0.
#  - Re-point the AutoFilter on column C (Scenario) from "National Trends" to
#    "Distributed Energy".
#  - Move the active selection to E12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row visibility -------------------------------------------------------
# Rows whose Scenario = "National Trends" become hidden.
$ws.Rows("2:6").Hidden = $true
$ws.Rows("15").Hidden = $true

# Rows whose Scenario = "Distributed Energy" become visible.
$ws.Rows("7:14").Hidden = $false
$ws.Rows("16:19").Hidden = $false

# --- New data row (20): NOS0 / hydrogen / Distributed Energy / 2040 / 500 -
$ws.Range("A20").Value = "NOS0"
$ws.Range("B20").Value = "hydrogen"
$ws.Range("C20").Value = "Distributed Energy"
$ws.Range("D20").Value = 2040
$ws.Range("E20").Value = 500
$ws.Rows("20").Hidden = $false

# --- AutoFilter: Scenario column now filters on "Distributed Energy" -----
$ws.Range("A1:E19").AutoFilter(3, @("Distributed Energy"), 7) | Out-Null

# --- Selection -------------------------------------------------------------
$ws.Range("E12").Select() | Out-Null
